$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Content Work for VA.gov Brand Consolidation: Vets.gov",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Content Work for VA.gov Brand Consolidation",
    2
)
